$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.429.78"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.721.65"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.95"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5321"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06724"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2672"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.07"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07726"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.523"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.957.19"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "1.711.07"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5878"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "0.0₅8242"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.31"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "27.436.49"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "225.53"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.683"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.54"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.068"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.35"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.703"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1214"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.279"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.39"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05398"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.297"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.496"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.437"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.634"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.873"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9592"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.389"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5913"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "1.154.99"
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01659"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.805"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8450"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.99"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.863.26"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.07"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4587"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.152"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05211"
$ws.Range("E51").Value = "  -0.57%  "
